$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.473.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.203.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0948"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +6.30%  "
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.531.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.208.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.485.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +20.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0760"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0300"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.201"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("E46").Value = "  +14.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.14%  "
